# Add a new "TaxableService" column (U) to the bulk invoice import template,
# mirroring the header formatting used by the existing last column (T).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell U1 holds the shared string "TaxableService".
# Setting .Value directly creates the shared-string entry, grows the
# sheet dimension/row spans, and inherits the row's header style (s="3"),
# matching how the rest of row 1 was produced.
$ws.Range("U1").Value = "TaxableService"

# Give the new column roughly the same "best fit" width behaviour as the
# other text-header columns (e.g. column T / InvoiceByHouseCode).
$ws.Columns.Item(21).ColumnWidth = 13.29

# Move the selection to the new last header cell, same as the author did
# after adding the column.
$ws.Range("U1").Select() | Out-Null

# Scroll the view right so the new column is visible (author's topLeftCell
# moved from column J to column L).
$excel.ActiveWindow.ScrollColumn = 12
